# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Pandaemonium Profits workbook
# as described by the commit diff (columns H-N across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2761.853
$ws.Range("I29").Value = 1380.6
$ws.Range("J29").Value = 3000
$ws.Range("K29").Value = 4141.799999999999
$ws.Range("L29").Value = 9000
$ws.Range("M29").Value = -3860.799999999999
$ws.Range("N29").Value = -9562

$ws.Range("H103").Value = 777.6
$ws.Range("I103").Value = 675
$ws.Range("J103").Value = 894.8570999999999
$ws.Range("K103").Value = 2025
$ws.Range("L103").Value = 2684.5713
$ws.Range("M103").Value = -1439
$ws.Range("N103").Value = -3856.5713

$ws.Range("H129").Value = 933.7231
$ws.Range("I129").Value = 323.6
$ws.Range("J129").Value = 984.56665
$ws.Range("K129").Value = 970.8000000000001
$ws.Range("L129").Value = 2953.69995
$ws.Range("M129").Value = 4029.2
$ws.Range("N129").Value = -12953.69995

$ws.Range("H137").Value = 743747.3
$ws.Range("I137").Value = 3986.6667
$ws.Range("J137").Value = 1113627.6
$ws.Range("K137").Value = 11960.0001
$ws.Range("L137").Value = 3340882.8
$ws.Range("M137").Value = -9410.000100000001
$ws.Range("N137").Value = -3345982.8

$ws.Range("H138").Value = 4363.2656
$ws.Range("I138").Value = 2058.4375
$ws.Range("J138").Value = 5131.5415
$ws.Range("K138").Value = 6175.3125
$ws.Range("L138").Value = 15394.6245
$ws.Range("M138").Value = -1035.3125
$ws.Range("N138").Value = -25674.6245

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22500.21
$ws.Range("I32").Value = 23264.959
$ws.Range("J32").Value = 10009.333
$ws.Range("K32").Value = 23264.959
$ws.Range("L32").Value = 10009.333
$ws.Range("M32").Value = -22977.959
$ws.Range("N32").Value = -10583.333

$ws.Range("H74").Value = 3670.4524
$ws.Range("I74").Value = 1301.4062
$ws.Range("J74").Value = 11251.4
$ws.Range("K74").Value = 1301.4062
$ws.Range("L74").Value = 11251.4
$ws.Range("M74").Value = -427.4061999999999
$ws.Range("N74").Value = -12999.4

$ws.Range("H77").Value = 3670.4524
$ws.Range("I77").Value = 1301.4062
$ws.Range("J77").Value = 11251.4
$ws.Range("K77").Value = 6507.030999999999
$ws.Range("L77").Value = 56257
$ws.Range("M77").Value = -2139.030999999999
$ws.Range("N77").Value = -64993

$ws.Range("H132").Value = 1847.791
$ws.Range("I132").Value = 1619.3208
$ws.Range("J132").Value = 2712.7144
$ws.Range("K132").Value = 4857.9624
$ws.Range("L132").Value = 8138.1432
$ws.Range("M132").Value = -2327.9624
$ws.Range("N132").Value = -13198.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 80780
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 80780
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 80780
$ws.Range("N53").Value = -81928

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 115
$ws.Range("I7").Value = 118
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 118
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = -5
$ws.Range("N7").Value = -326

$ws.Range("H31").Value = 732098.7
$ws.Range("I31").Value = 14359.4
$ws.Range("J31").Value = 1031156.75
$ws.Range("K31").Value = 14359.4
$ws.Range("L31").Value = 1031156.75
$ws.Range("M31").Value = -14064.4
$ws.Range("N31").Value = -1031746.75

$ws.Range("H34").Value = 732098.7
$ws.Range("I34").Value = 14359.4
$ws.Range("J34").Value = 1031156.75
$ws.Range("K34").Value = 14359.4
$ws.Range("L34").Value = 1031156.75
$ws.Range("M34").Value = -14157.4
$ws.Range("N34").Value = -1031560.75

$ws.Range("H50").Value = 20307.092
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 20307.092
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 20307.092
$ws.Range("N50").Value = -21557.092

$ws.Range("H51").Value = 22373.092
$ws.Range("I51").Value = 7090
$ws.Range("J51").Value = 23901.4
$ws.Range("K51").Value = 7090
$ws.Range("L51").Value = 23901.4
$ws.Range("M51").Value = -6354
$ws.Range("N51").Value = -25373.4

$ws.Range("H58").Value = 2168714.5
$ws.Range("I58").Value = 3031702
$ws.Range("J58").Value = 11246.167
$ws.Range("K58").Value = 3031702
$ws.Range("L58").Value = 11246.167
$ws.Range("M58").Value = -3031499
$ws.Range("N58").Value = -11652.167

$ws.Range("H59").Value = 25112.611
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 25112.611
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 25112.611
$ws.Range("M59").Value = $null
$ws.Range("N59").Value = -27402.611

$ws.Range("H60").Value = 24266.666
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 24266.666
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 24266.666
$ws.Range("N60").Value = -25288.666

$ws.Range("H61").Value = 22373.092
$ws.Range("I61").Value = 7090
$ws.Range("J61").Value = 23901.4
$ws.Range("K61").Value = 7090
$ws.Range("L61").Value = 23901.4
$ws.Range("M61").Value = -6742
$ws.Range("N61").Value = -24597.4

$ws.Range("H68").Value = 40295
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 40295
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 40295
$ws.Range("N68").Value = -41793

$ws.Range("H71").Value = 40295
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 40295
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 120885
$ws.Range("N71").Value = -128373

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null

$ws.Range("H94").Value = 382.84616
$ws.Range("I94").Value = 200
$ws.Range("J94").Value = 398.08334
$ws.Range("K94").Value = 200
$ws.Range("L94").Value = 398.08334
$ws.Range("M94").Value = 251
$ws.Range("N94").Value = -1300.08334

$ws.Range("H136").Value = 2168714.5
$ws.Range("I136").Value = 3031702
$ws.Range("J136").Value = 11246.167
$ws.Range("K136").Value = 9095106
$ws.Range("L136").Value = 33738.501
$ws.Range("M136").Value = -9092556
$ws.Range("N136").Value = -38838.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 32.942856
$ws.Range("I2").Value = 11.75
$ws.Range("J2").Value = 35.677418
$ws.Range("K2").Value = 70.5
$ws.Range("L2").Value = 214.064508
$ws.Range("M2").Value = 42.5
$ws.Range("N2").Value = -440.064508

$ws.Range("H5").Value = 16677138
$ws.Range("I5").Value = 501.6
$ws.Range("J5").Value = 33353776
$ws.Range("K5").Value = 1504.8
$ws.Range("L5").Value = 100061328
$ws.Range("M5").Value = -1392.8
$ws.Range("N5").Value = -100061552

$ws.Range("H63").Value = 2353
$ws.Range("I63").Value = 1520.6666
$ws.Range("J63").Value = 4850
$ws.Range("K63").Value = 4561.9998
$ws.Range("L63").Value = 14550
$ws.Range("M63").Value = -3812.9998
$ws.Range("N63").Value = -16048

$ws.Range("H66").Value = 2353
$ws.Range("I66").Value = 1520.6666
$ws.Range("J66").Value = 4850
$ws.Range("K66").Value = 13685.9994
$ws.Range("L66").Value = 43650
$ws.Range("M66").Value = -9941.999400000001
$ws.Range("N66").Value = -51138

$ws.Range("H75").Value = 3665.7
$ws.Range("I75").Value = 1978.25
$ws.Range("J75").Value = 4790.6665
$ws.Range("K75").Value = 5934.75
$ws.Range("L75").Value = 14371.9995
$ws.Range("M75").Value = -4936.75
$ws.Range("N75").Value = -16367.9995

$ws.Range("H78").Value = 3665.7
$ws.Range("I78").Value = 1978.25
$ws.Range("J78").Value = 4790.6665
$ws.Range("K78").Value = 17804.25
$ws.Range("L78").Value = 43115.9985
$ws.Range("M78").Value = -12812.25
$ws.Range("N78").Value = -53099.9985

$ws.Range("H81").Value = 4673.636
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 4941
$ws.Range("K81").Value = 6000
$ws.Range("L81").Value = 14823
$ws.Range("M81").Value = -4877
$ws.Range("N81").Value = -17069

$ws.Range("H84").Value = 4673.636
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 4941
$ws.Range("K84").Value = 18000
$ws.Range("L84").Value = 44469
$ws.Range("M84").Value = -12384
$ws.Range("N84").Value = -55701

$ws.Range("H93").Value = 5005.4
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 5005.4
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 15016.2
$ws.Range("N93").Value = -18760.2

$ws.Range("H96").Value = 4336.1113
$ws.Range("I96").Value = 1025
$ws.Range("J96").Value = 4750
$ws.Range("K96").Value = 3075
$ws.Range("L96").Value = 14250
$ws.Range("M96").Value = -1016
$ws.Range("N96").Value = -18368

$ws.Range("H105").Value = 80000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 80000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 240000
$ws.Range("N105").Value = -245242

$ws.Range("H113").Value = 568.77
$ws.Range("I113").Value = 636.6923
$ws.Range("J113").Value = 495.1875
$ws.Range("K113").Value = 1910.0769
$ws.Range("L113").Value = 1485.5625
$ws.Range("M113").Value = 259.9231
$ws.Range("N113").Value = -5825.5625

$ws.Range("H123").Value = 2200
$ws.Range("I123").Value = 400
$ws.Range("J123").Value = 4000
$ws.Range("K123").Value = 1200
$ws.Range("L123").Value = 12000
$ws.Range("M123").Value = 1250
$ws.Range("N123").Value = -16900

$ws.Range("H125").Value = 2571.4
$ws.Range("I125").Value = 2571.4
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 7714.200000000001
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -2794.200000000001

$ws.Range("H131").Value = 41459.332
$ws.Range("I131").Value = 2243.75
$ws.Range("J131").Value = 61067.125
$ws.Range("K131").Value = 6731.25
$ws.Range("L131").Value = 183201.375
$ws.Range("M131").Value = -1691.25
$ws.Range("N131").Value = -193281.375

$ws.Range("H135").Value = 16677138
$ws.Range("I135").Value = 501.6
$ws.Range("J135").Value = 33353776
$ws.Range("K135").Value = 4514.400000000001
$ws.Range("L135").Value = 300183984
$ws.Range("M135").Value = -1979.400000000001
$ws.Range("N135").Value = -300189054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 55980
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 55980
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 55980
$ws.Range("N123").Value = -65780

$ws.Range("H132").Value = 7242.067
$ws.Range("I132").Value = 8432.041999999999
$ws.Range("J132").Value = 2482.1667
$ws.Range("K132").Value = 25296.126
$ws.Range("L132").Value = 7446.500100000001
$ws.Range("M132").Value = -22766.126
$ws.Range("N132").Value = -12506.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4999.5
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 4999.5
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 4999.5
$ws.Range("M14").Value = $null
$ws.Range("N14").Value = -5335.5

$ws.Range("H64").Value = 38114
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 38114
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 38114
$ws.Range("N64").Value = -38610

$ws.Range("H67").Value = 38114
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 38114
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 38114
$ws.Range("N67").Value = -39830

$ws.Range("H127").Value = 62714.5
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 62714.5
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 62714.5
$ws.Range("N127").Value = -72634.5
